$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = 2305722
$ws.Range("B2").Value = 285477

# Row 3 updates
$ws.Range("A3").Value = 2312187
$ws.Range("B3").Value = 27388

# Fill in rows 4-7 with new data.
# Copy the C:E formatting/values from row 2 (style index 14) onto rows 4-7
# so the number format/fill style matches the target layout, then set A/B.
$srcCE = $ws.Range("C2:E2")
$srcCE.Copy($ws.Range("C4:E4"))
$srcCE.Copy($ws.Range("C5:E5"))
$srcCE.Copy($ws.Range("C6:E6"))
$srcCE.Copy($ws.Range("C7:E7"))

$ws.Range("A4").Value = 2312608
$ws.Range("B4").Value = 10723

$ws.Range("A5").Value = 2316491
$ws.Range("B5").Value = 23289

$ws.Range("A6").Value = 2319581
$ws.Range("B6").Value = 16151

$ws.Range("A7").Value = 2319636
$ws.Range("B7").Value = 73757

# Update the active selection to B7, matching the saved sheet view state.
$ws.Range("B7").Select()
